$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan")

$xlPasteFormats = -4122

# C3 formula becomes an absolute reference to A32
$ws.Range("C3").Formula = "=`$A32"

# D3 gets the same formula, formatted like C3 (numFmtId 0 / right aligned)
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial($xlPasteFormats)
$ws.Range("D3").Formula = "=`$A32"

# New time-of-day entries recorded in column D, formatted like their column C counterparts
$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial($xlPasteFormats)
$ws.Range("D4").Value = 0.54652777777777783

$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial($xlPasteFormats)
$ws.Range("D5").Value = 0.60902777777777783

$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial($xlPasteFormats)
$ws.Range("D6").Value = 0.65138888888888891

$ws.Range("C4").Copy()
$ws.Range("D7").PasteSpecial($xlPasteFormats)
$ws.Range("D7").Value = 0.66805555555555562

$ws.Range("C4").Copy()
$ws.Range("D8").PasteSpecial($xlPasteFormats)
$ws.Range("D8").Value = 0.84375

$ws.Range("C21").Copy()
$ws.Range("D21").PasteSpecial($xlPasteFormats)
$ws.Range("D21").Value = 0.88888888888888884

$excel.CutCopyMode = $false

# Move the active selection to D21, matching the saved view state
$ws.Range("D21").Select()
